# Auto-generated Excel COM-interop script to apply the diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 11217
$ws.Range("F4").Value = 281
$ws.Range("F5").Value = 1260
$ws.Range("F6").Value = 1138
$ws.Range("F7").Value = 876
$ws.Range("F8").Value = 299
$ws.Range("F10").Value = 1200
$ws.Range("F11").Value = 161
$ws.Range("F12").Value = 932
$ws.Range("F13").Value = 2177
$ws.Range("F14").Value = 27
$ws.Range("F15").Value = 1074
$ws.Range("F16").Value = 858
$ws.Range("F17").Value = 570
$ws.Range("F18").Value = 840
$ws.Range("F19").Value = 974
$ws.Range("F21").Value = 275
$ws.Range("F22").Value = 96
$ws.Range("F23").Value = 670
$ws.Range("F24").Value = 696
$ws.Range("F26").Value = 383
$ws.Range("F27").Value = 1038
$ws.Range("F28").Value = 56
$ws.Range("F30").Value = 520
$ws.Range("F33").Value = 260
$ws.Range("F34").Value = 611
$ws.Range("F35").Value = 2275
$ws.Range("F36").Value = 417
$ws.Range("F37").Value = 59
$ws.Range("F38").Value = 1487
$ws.Range("F39").Value = 418
$ws.Range("F41").Value = 61
$ws.Range("F43").Value = 51
$ws.Range("F45").Value = 92
$ws.Range("F47").Value = 60
$ws.Range("F48").Value = 14
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 209
$ws.Range("F11").Value = 92
$ws.Range("C17").Value = "广州·音阅派国漫演唱会-《狐妖小红娘》《一人之下》领衔国漫原声音乐现场"
$ws.Range("F17").Value = 10
$ws.Range("F18").Value = 20
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2209
$ws.Range("F3").Value = 664
$ws.Range("F4").Value = 613
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2209
$ws.Range("F4").Value = 281
$ws.Range("F5").Value = 1260
$ws.Range("F6").Value = 613
$ws.Range("F7").Value = 1138
$ws.Range("F8").Value = 876
$ws.Range("F9").Value = 209
$ws.Range("F10").Value = 299
$ws.Range("F11").Value = 1200
$ws.Range("F14").Value = 932
$ws.Range("F15").Value = 2178
$ws.Range("F16").Value = 27
$ws.Range("F17").Value = 1074
$ws.Range("F18").Value = 858
$ws.Range("F19").Value = 570
$ws.Range("F20").Value = 840
$ws.Range("F21").Value = 974
$ws.Range("F22").Value = 275
$ws.Range("F24").Value = 96
$ws.Range("F25").Value = 670
$ws.Range("F26").Value = 696
$ws.Range("F28").Value = 383
$ws.Range("F29").Value = 1038
$ws.Range("F30").Value = 56
$ws.Range("F32").Value = 520
$ws.Range("F35").Value = 260
$ws.Range("F36").Value = 2275
$ws.Range("F38").Value = 417
$ws.Range("F39").Value = 59
$ws.Range("F40").Value = 1487
$ws.Range("F41").Value = 418
$ws.Range("F44").Value = 51
$ws.Range("F45").Value = 92
$ws.Range("F47").Value = 60
